$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update "VALOR MORA" total value (E11)
$ws.Range("E11").Value = 409824

# 2) Update "Cant. Periodos" value (F13)
$ws.Range("F13").Value = 6

# 3) Duplicate the last worker-period row (row 20) into a new row for period 2509.
#    First, snapshot row 20 (values + formats) so we can clone it into the newly
#    inserted row.
$ws.Range("B20:J20").Copy()

# Insert a new blank row at 21; this pushes the old rows 21-26 down to 22-27.
$ws.Rows("21:21").Insert()

# Paste the captured values then formats from the old row 20 into the new row 21,
# so row 21 becomes an exact clone of what row 20 used to be (same "last row"
# bottom-border styling and same worker info/amounts).
$ws.Range("B21:J21").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)  # xlPasteFormats

# Update the new row's period to the new one (2509).
$ws.Range("E21").Value = "2509"

# Row 20 is no longer the last data row, so restyle it like the other interior
# rows (copy formatting from row 19, which still has the "middle row" style).
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# The "Periodo Mora" column is now center-aligned across all worker rows.
$ws.Range("E16:E21").HorizontalAlignment = -4108  # xlCenter
